$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before U. This shifts the old "Original - Coretop"
# header/values from U to V, and the old "BAYMAG - Coretop" header/values
# from V to W, leaving a blank column U ready for the new data.
$ws.Range("U1").EntireColumn.Insert()

# Header for the newly inserted column U.
$ws.Range("U1").Value = "MgCa Coretop modelled temperature"

# Value for the newly inserted column in the data row.
$ws.Range("U2").Value = 28.7843

# Updated values for existing columns in the data row.
$ws.Range("M2").Value = 29.37598672
$ws.Range("N2").Value = 30.8326333333333
$ws.Range("R2").Value = 29.09
$ws.Range("S2").Value = 0.282595301542969
$ws.Range("T2").Value = 1.739241914876267
$ws.Range("V2").Value = 0.5916867199999984
$ws.Range("W2").Value = 2.048333333333296
